$d = $word.ActiveDocument

$replacements = @(
    @("390÷8=48, 6", "184÷7=26, 2"),
    @("660÷2=330, 0", "828÷8=103, 4"),
    @("929÷8=116, 1", "390÷4=97, 2"),
    @("291÷2=145, 1", "972÷7=138, 6"),
    @("567÷8=70, 7", "831÷9=92, 3"),
    @("687÷3=229, 0", "176÷7=25, 1"),
    @("861÷4=215, 1", "902÷9=100, 2"),
    @("567÷4=141, 3", "407÷2=203, 1"),
    @("383÷7=54, 5", "465÷5=93, 0"),
    @("869÷3=289, 2", "579÷2=289, 1"),
    @("420÷7=60, 0", "748÷8=93, 4"),
    @("406÷8=50, 6", "770÷2=385, 0"),
    @("555÷9=61, 6", "153÷8=19, 1"),
    @("559÷4=139, 3", "131÷2=65, 1"),
    @("539÷2=269, 1", "450÷4=112, 2"),
    @("677÷9=75, 2", "359÷8=44, 7"),
    @("709÷6=118, 1", "492÷4=123, 0"),
    @("393÷6=65, 3", "615÷2=307, 1"),
    @("391÷9=43, 4", "894÷7=127, 5"),
    @("489÷2=244, 1", "244÷8=30, 4"),
    @("836÷4=209, 0", "500÷6=83, 2"),
    @("128÷7=18, 2", "644÷8=80, 4"),
    @("547÷5=109, 2", "359÷6=59, 5"),
    @("665÷4=166, 1", "392÷9=43, 5"),
    @("950÷5=190, 0", "450÷4=112, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
